$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells in column D keep their exact literal formatting
# (European-style thousand separators, fixed decimal places, etc.) by forcing
# the Text number format before assigning, so Excel does not coerce the
# strings into numeric values.

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '74.841.38'
$ws.Cells.Item(2, 5).Value = '  +0.87%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.813.38'
$ws.Cells.Item(3, 5).Value = '  +6.77%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '187.04'
$ws.Cells.Item(5, 5).Value = '  +0.53%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '593.84'
$ws.Cells.Item(6, 5).Value = '  +1.89%  '
$ws.Cells.Item(7, 5).Value = '  +0.00%  '
$ws.Cells.Item(8, 5).Value = '  +2.58%  '
$ws.Cells.Item(9, 5).Value = '  -4.97%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '2.812.11'
$ws.Cells.Item(10, 5).Value = '  +6.79%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.161'
$ws.Cells.Item(11, 5).Value = '  -1.27%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.370'
$ws.Cells.Item(12, 5).Value = '  +3.45%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '4.88'
$ws.Cells.Item(13, 5).Value = '  +2.39%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '3.332.64'
$ws.Cells.Item(14, 5).Value = '  +6.92%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '74.697.51'
$ws.Cells.Item(15, 5).Value = '  +1.20%  '
$ws.Cells.Item(16, 5).Value = '  -1.46%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '26.74'
$ws.Cells.Item(17, 5).Value = '  +1.64%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.814.94'
$ws.Cells.Item(18, 5).Value = '  +7.18%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '8.94'
$ws.Cells.Item(19, 5).Value = '  -1.70%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '12.28'
$ws.Cells.Item(20, 5).Value = '  +3.86%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '377.10'
$ws.Cells.Item(21, 5).Value = '  +1.18%  '
$ws.Cells.Item(22, 5).Value = '  -1.70%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '4.07'
$ws.Cells.Item(23, 5).Value = '  -0.62%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '0.999'
$ws.Cells.Item(24, 5).Value = '  -0.22%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '70.81'
$ws.Cells.Item(25, 5).Value = '  +0.92%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '9.88'
$ws.Cells.Item(26, 5).Value = '  +5.41%  '
$ws.Cells.Item(27, 5).Value = '  +7.06%  '
$ws.Cells.Item(28, 5).Value = '  -0.20%  '
$ws.Cells.Item(29, 5).Value = '  +9.16%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.996'
$ws.Cells.Item(30, 5).Value = '  -1.33%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '517.95'
$ws.Cells.Item(31, 5).Value = '  -1.62%  '
$ws.Cells.Item(32, 5).Value = '  -0.54%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '7.74'
$ws.Cells.Item(33, 5).Value = '  +0.19%  '
$ws.Cells.Item(34, 5).Value = '  +2.35%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.999'
$ws.Cells.Item(35, 5).Value = '  -0.03%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '163.34'
$ws.Cells.Item(36, 5).Value = '  +0.12%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '19.91'
$ws.Cells.Item(37, 5).Value = '  +3.81%  '
$ws.Cells.Item(38, 5).Value = '  -1.07%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '19.36'
$ws.Cells.Item(39, 5).Value = '  +0.41%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '185.84'
$ws.Cells.Item(40, 5).Value = '  +15.60%  '
$ws.Cells.Item(41, 5).Value = '  +0.02%  '
$ws.Cells.Item(42, 5).Value = '  +3.95%  '
$ws.Cells.Item(43, 5).Value = '  +1.50%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.66'
$ws.Cells.Item(44, 5).Value = '  -0.65%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '39.92'
$ws.Cells.Item(46, 5).Value = '  +2.49%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0856'
$ws.Cells.Item(47, 5).Value = '  +0.04%  '
$ws.Cells.Item(48, 5).Value = '  -2.86%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.575'
$ws.Cells.Item(49, 5).Value = '  +8.47%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '3.70'
$ws.Cells.Item(50, 5).Value = '  +2.21%  '
$ws.Cells.Item(51, 5).Value = '  +7.88%  '
